$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("model_3_5_20", 0.4367514834460438, 0.2952304384545075, 0.4522679027324807, 0.3847582114666007, 0.6233501434326172, 0.7578065395355225, 0.8220329284667969, 0.7880306839942932),
    @("model_3_5_24", 0.4370858593187242, 0.2937061881377719, 0.4544790886311105, 0.3853007057528003, 0.6229801177978516, 0.7594454884529114, 0.8187143802642822, 0.7873358130455017),
    @("model_3_5_23", 0.4373199646620953, 0.2931565804556433, 0.4562952523952334, 0.3860578209621109, 0.6227210164070129, 0.7600364089012146, 0.8159886598587036, 0.7863660454750061),
    @("model_3_5_22", 0.4374164490946596, 0.2939635143771404, 0.4562502608915323, 0.3863902349954188, 0.6226142644882202, 0.759168803691864, 0.8160561919212341, 0.7859402894973755),
    @("model_3_5_21", 0.4374256110570883, 0.2945403851669307, 0.4560084516119806, 0.3865139670987241, 0.6226040720939636, 0.7585484981536865, 0.8164191246032715, 0.7857817411422729),
    @("model_3_5_17", 0.4382116227194593, 0.2983485763218057, 0.458718350299772, 0.3897008364085585, 0.6217342019081116, 0.7544537782669067, 0.8123520612716675, 0.7816998958587646),
    @("model_3_5_18", 0.4382753309025219, 0.2980978498547525, 0.4591029433436402, 0.3898012238899734, 0.6216636896133423, 0.7547233104705811, 0.8117748498916626, 0.7815713286399841),
    @("model_3_5_19", 0.4383506705722905, 0.2967174704321972, 0.4601408908350862, 0.3897601083036498, 0.6215803027153015, 0.7562076449394226, 0.8102171421051025, 0.7816239595413208),
    @("model_3_5_15", 0.4396217544343949, 0.298174623089792, 0.4666886704678475, 0.3940182056207997, 0.6201735734939575, 0.7546408176422119, 0.8003903031349182, 0.7761700749397278),
    @("model_3_5_16", 0.4403149357244391, 0.2973344490877726, 0.4704061084316988, 0.3956946236713801, 0.6194064021110535, 0.7555441856384277, 0.7948111891746521, 0.7740227580070496),
    @("model_3_5_14", 0.440722030933982, 0.2953203791718995, 0.4756056848553092, 0.397666569250931, 0.6189559102058411, 0.757709801197052, 0.7870076894760132, 0.7714970111846924),
    @("model_3_5_13", 0.4419984626351456, 0.2970098590517163, 0.4821837401722917, 0.4020445007804886, 0.6175432801246643, 0.7558932304382324, 0.7771353721618652, 0.765889585018158),
    @("model_3_5_11", 0.4435004721709457, 0.2974970902196133, 0.491585452446322, 0.4074452620957844, 0.6158810257911682, 0.7553693056106567, 0.7630254030227661, 0.7589719891548157),
    @("model_3_5_12", 0.4439535414429654, 0.2975622333105741, 0.4935529839927255, 0.4085590049878508, 0.6153796315193176, 0.7552992105484009, 0.7600725293159485, 0.7575454711914062),
    @("model_3_5_1", 0.4463599943165895, 0.3095808913288109, 0.591176650730421, 0.4677298307450463, 0.6127163767814636, 0.7423761487007141, 0.6135594844818115, 0.6817566156387329),
    @("model_3_5_10", 0.4477237636180231, 0.3005763718055854, 0.5134206932999195, 0.4208536188409819, 0.6112070083618164, 0.7520582675933838, 0.7302551865577698, 0.741797924041748),
    @("model_3_5_9", 0.449023496146748, 0.3009302566003308, 0.5209163528935264, 0.4251438834101413, 0.6097686886787415, 0.7516777515411377, 0.7190057039260864, 0.7363027930259705),
    @("model_3_5_2", 0.4498441431535528, 0.3111728333808697, 0.5971446758045701, 0.4717282067918867, 0.6088604331016541, 0.7406643629074097, 0.6046026945114136, 0.6766353249549866),
    @("model_3_5_8", 0.4498644550309888, 0.3049715012274276, 0.5248161192806875, 0.4290903016276619, 0.608837902545929, 0.7473324537277222, 0.7131530046463013, 0.7312480211257935),
    @("model_3_5_5", 0.4500445676355206, 0.3113079611938683, 0.5512475901059746, 0.4464807723038128, 0.6086386442184448, 0.740519106388092, 0.673484742641449, 0.7089734673500061),
    @("model_3_5_4", 0.4510406980852047, 0.3109408997959683, 0.5630057252118896, 0.4528009016469473, 0.6075361967086792, 0.7409137487411499, 0.6558383107185364, 0.700878381729126),
    @("model_3_5_7", 0.4511701404808983, 0.308875287911469, 0.5344522999174089, 0.4361384485151414, 0.6073929071426392, 0.7431348562240601, 0.6986910700798035, 0.7222204208374023),
    @("model_3_5_6", 0.4528099687495081, 0.3094614540586423, 0.5492429062804431, 0.4445546658385123, 0.6055780649185181, 0.7425045967102051, 0.6764934659004211, 0.7114405632019043),
    @("model_3_5_3", 0.4530954830859569, 0.3130454765828736, 0.5912187461698611, 0.4692929181181529, 0.6052622199058533, 0.738650918006897, 0.6134963035583496, 0.6797546148300171),
    @("model_3_5_0", 0.4601610986532527, 0.3214512321410198, 0.7075328688409075, 0.5371640293382054, 0.597442626953125, 0.7296125292778015, 0.4389328360557556, 0.5928220748901367)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    $ws.Cells.Item($r, 1).Value = $rowVals[0]
    for ($c = 1; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = [double]$rowVals[$c]
    }
}